$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the Hawman et al. mouse model row (row 5), which no longer has
# sequenced GenBank accessions for the mouse-adapted CCHFV model.
$ws.Rows.Item(5).Delete()
